$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript digit characters used in some token price cells (e.g. PEPE)
$sub3 = [char]0x2083
$sub6 = [char]0x2086

# Force the Price/Volume cells being updated to remain plain text so Excel
# does not reinterpret values like "1.00" or "0.0000144" as numbers/dates.
$ws.Range("D2,E2,D3,E3,D4,E4,D5,E5,D6,E6,D8,E8,D9,E9,D10,E10,D11,E11,E12,D13,E13,D14,E14,D15,E15,D16,E16,D17,E17,D18,E18,D19,E19,D20,E20,D21,E21,D22,E22,D23,E23,D24,E24,E25,D26,E26,D27,E27,D28,E28,D29,E29,D30,E30,D31,E31,D32,E32,D33,E33,D34,E34,D35,E35,D36,E36,D37,E37,D38,E38,D39,E39,D40,E40,D41,E41,D42,E42,D43,E43,D44,E44,D45,E45,D46,E46,D47,E47,D48,E48,D49,E49,D50,E50,E51").NumberFormat = "@"

$ws.Range("D2").Value = '64.041.45'
$ws.Range("E2").Value = '  +3.22%  '
$ws.Range("D3").Value = '2.539.21'
$ws.Range("E3").Value = '  +5.76%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '572.76'
$ws.Range("E5").Value = '  +2.27%  '
$ws.Range("D6").Value = '146.11'
$ws.Range("E6").Value = '  +5.83%  '
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  +0.82%  '
$ws.Range("D9").Value = '2.537.57'
$ws.Range("E9").Value = '  +5.75%  '
$ws.Range("D10").Value = '0.107'
$ws.Range("E10").Value = '  +2.56%  '
$ws.Range("D11").Value = '5.79'
$ws.Range("E11").Value = '  +1.32%  '
$ws.Range("E12").Value = '  +1.76%  '
$ws.Range("D13").Value = '0.360'
$ws.Range("E13").Value = '  +3.22%  '
$ws.Range("D14").Value = '28.07'
$ws.Range("E14").Value = '  +9.64%  '
$ws.Range("D15").Value = '2.983.74'
$ws.Range("E15").Value = '  +5.42%  '
$ws.Range("D16").Value = '63.815.03'
$ws.Range("E16").Value = '  +2.87%  '
$ws.Range("D17").Value = '0.0000144'
$ws.Range("E17").Value = '  +4.41%  '
$ws.Range("D18").Value = '2.520.30'
$ws.Range("E18").Value = '  +5.04%  '
$ws.Range("D19").Value = '11.49'
$ws.Range("E19").Value = '  +4.51%  '
$ws.Range("D20").Value = '344.56'
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").Value = '4.36'
$ws.Range("E21").Value = '  +3.34%  '
$ws.Range("D22").Value = '6.91'
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("D23").Value = '1.01'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("D24").Value = '66.10'
$ws.Range("E24").Value = '  +1.73%  '
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").Value = '1.57'
$ws.Range("E26").Value = '  +5.37%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '8.25'
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("B29").Value = 'SuiNetwork'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D29").Value = '1.43'
$ws.Range("E29").Value = '  +4.38%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0' + [string]$sub3 + '0827'
$ws.Range("E30").Value = '  +7.50%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '1.88'
$ws.Range("E31").Value = '  +4.46%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").Value = '6.80'
$ws.Range("E32").Value = '  +6.93%  '
$ws.Range("D33").Value = '176.44'
$ws.Range("E33").Value = '  +2.86%  '
$ws.Range("D34").Value = '1.55'
$ws.Range("E34").Value = '  +10.17%  '
$ws.Range("D35").Value = '410.42'
$ws.Range("E35").Value = '  +14.35%  '
$ws.Range("D36").Value = '0.403'
$ws.Range("E36").Value = '  +2.75%  '
$ws.Range("D37").Value = '19.19'
$ws.Range("E37").Value = '  +3.81%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '4.44'
$ws.Range("E38").Value = '  -2.00%  '
$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").Value = '1.76'
$ws.Range("E40").Value = '  +5.46%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").Value = '40.77'
$ws.Range("E42").Value = '  +4.65%  '
$ws.Range("D43").Value = '152.81'
$ws.Range("E43").Value = '  +6.68%  '
$ws.Range("D44").Value = '3.79'
$ws.Range("E44").Value = '  +3.65%  '
$ws.Range("D45").Value = '21.08'
$ws.Range("E45").Value = '  +3.58%  '
$ws.Range("D46").Value = '0.615'
$ws.Range("E46").Value = '  +5.58%  '
$ws.Range("D47").Value = '0.0968'
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("D48").Value = '0.0529'
$ws.Range("E48").Value = '  +2.12%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '18.99'
$ws.Range("E49").Value = '  +6.87%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '0.0231'
$ws.Range("E50").Value = '  +4.48%  '
$ws.Range("E51").Value = '  +5.85%  '
